# fix bugs, update isolate column comment if the isolate is not a clinical isolate.
# Capitalize the first letter of the "Specimen" (column AB) value for every
# data row, except where the value is "NA" (not applicable / unknown).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 28)  # column AB = 28
    $text = $cell.Text

    if ($text -and -not $text.Equals("NA")) {
        $first = $text.Substring(0, 1).ToUpper()
        $rest = $text.Substring(1)
        $newValue = $first + $rest
        $cell.Value = $newValue
    }
}
